$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 220, shifting existing rows 220:258 down to 221:259
$ws.Rows.Item(220).Insert()

# Populate the newly inserted row 220 with the new weekly record
$ws.Cells.Item(220, 1).Value = 8
$ws.Cells.Item(220, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(220, 3).Value = "Coquimbo"
$ws.Cells.Item(220, 4).Value = 45218
$ws.Cells.Item(220, 5).Value = 4
$ws.Cells.Item(220, 6).Value = 100112044
$ws.Cells.Item(220, 7).Value = "Perejil"
$ws.Cells.Item(220, 8).Value = "Sin especificar"
$ws.Cells.Item(220, 9).Value = "Primera"
$ws.Cells.Item(220, 10).Value = 2000
$ws.Cells.Item(220, 11).Value = 1500
$ws.Cells.Item(220, 12).Value = 2000
$ws.Cells.Item(220, 13).Value = 1750
$ws.Cells.Item(220, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(220, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(220, 16).Value = 1167
$ws.Cells.Item(220, 17).Value = 1.5
$ws.Cells.Item(220, 18).Value = "Hortaliza"
